$p = $ppt.ActivePresentation

$oldDate = "2023/12/10"
$newDate = "2023/12/21"
$ppPlaceholderDate = 16

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePh = $true
            }
        } catch {
        }
        if ($isDatePh -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Update the cached datetimeFigureOut field text on the slide master ...
$master = $p.Slides.Item(1).Master
Update-DateShape $master.Shapes

# ... and on every slide layout that belongs to it.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShape $layout.Shapes
}

# Rename the title on slide 1 from "挑戰模式" to "重新開始".
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "挑戰模式") {
        $shp.TextFrame.TextRange.Text = "重新開始"
    }
}
